# The "ExtractionType" (I) and "SamplePortion" (J) columns are being
# removed from the chromatography data_info template; remaining columns
# (in particular "Comment", previously K) shift left to fill the gap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1:J3").Delete()
